$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.836.41'
$ws.Range('E2').Value = '  -1.31%  '
$ws.Range('D3').Value = '3.502.28'
$ws.Range('E3').Value = '  -1.62%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '608.06'
$ws.Range('E5').Value = '  +3.85%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '191.62'
$ws.Range('E6').Value = '  +1.77%  '
$ws.Range('E7').Value = '  +0.65%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  -0.21%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.665'
$ws.Range('E10').Value = '  +3.22%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '53.45'
$ws.Range('E11').Value = '  -0.86%  '
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '9.61'
$ws.Range('E13').Value = '  +2.28%  '
$ws.Range('D14').Value = '4.061.45'
$ws.Range('E14').Value = '  -1.56%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '621.94'
$ws.Range('E15').Value = '  +9.95%  '
$ws.Range('D16').Value = '69.911.75'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '12.71'
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('D19').Value = '3.506.48'
$ws.Range('E19').Value = '  -2.17%  '
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '17.73'
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('E23').Value = '  +12.91%  '
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.00'
$ws.Range('E25').Value = '  +2.61%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.04'
$ws.Range('E26').Value = '  +4.84%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.98'
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('E28').Value = '  +5.54%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '34.16'
$ws.Range('E29').Value = '  +5.58%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.06'
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '12.59'
$ws.Range('E31').Value = '  +3.52%  '
$ws.Range('E32').Value = '  +4.84%  '
$ws.Range('E33').Value = '  +0.27%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '64.29'
$ws.Range('E34').Value = '  +1.91%  '
$ws.Range('D35').Value = '3.710.51'
$ws.Range('E35').Value = '  +1.95%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.10'
$ws.Range('E36').Value = '  -4.24%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '519.14'
$ws.Range('E38').Value = '  -1.74%  '
$ws.Range('D39').Value = '0.0₃0794'
$ws.Range('E39').Value = '  +1.16%  '
$ws.Range('E40').Value = '  -3.71%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '36.78'
$ws.Range('E41').Value = '  -3.53%  '
$ws.Range('E42').Value = '  +1.25%  '
$ws.Range('E43').Value = '  -0.66%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0462'
$ws.Range('E44').Value = '  +1.03%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.85'
$ws.Range('E45').Value = '  -2.10%  '
$ws.Range('E46').Value = '  +2.50%  '
$ws.Range('E47').Value = '  -3.68%  '
$ws.Range('E48').Value = '  +0.39%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.74'
$ws.Range('E49').Value = '  -4.78%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '132.33'
$ws.Range('E50').Value = '  -2.07%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.30'
$ws.Range('E51').Value = '  +10.94%  '
